$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "NAZARENA RAOS"
$ws.Range("B14").Value = "Daniele  Dalbosco | iMontagna"
$ws.Range("C14").Value = "ANDREA ASTE | Pinguini Trentini"
$ws.Range("D14").Value = "Matteo Zanlucchi | SBARX"
$ws.Range("E14").Value = "Federico  Mortillaro | Clitoriders"
$ws.Range("F14").Value = "Stefano Mattioli | MAI UNA GIOIA"
